# warning-report-template.docx:
#  1. Add a "Failed Expressions" / {{FailedExpressionCount}} row to the
#     summary table (end of table 1).
#  2. Add a whole new "{{#if HasFailedExpressions}} ... {{/if}}" section
#     (mirroring the existing Missing Variables / Missing Collections /
#     Null Collections sections) right after the HasNullCollections
#     section and before the trailing "End of Warning Report" paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append a row to the summary table (the first table in the doc).
# ---------------------------------------------------------------------
$summaryTable = $d.Tables(1)
$newRow = $summaryTable.Rows.Add()
$newRow.Cells(1).Split(1, 2)
$newRow.Cells(1).Range.Text = "Failed Expressions"
$newRow.Cells(2).Range.Text = "{{FailedExpressionCount}}"

# ---------------------------------------------------------------------
# 2. Insert the new "Failed Expressions" detail section.
#
# It goes right after the paragraph that closes the HasNullCollections
# section (the *last* "{{/if}}" paragraph in the document) and before
# the blank paragraph that precedes "End of Warning Report".
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$paraCount = $paras.Count
$endReportIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $paras.Item($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq "End of Warning Report") {
        $endReportIndex = $i
        break
    }
}

# The blank paragraph right before "End of Warning Report" - insert at its
# Start (== the End of the preceding "{{/if}}" paragraph) so the new
# content lands as its own paragraph(s) ahead of that blank separator.
$blankPara = $paras.Item($endReportIndex - 1)
$insertBefore = $d.Range($blankPara.Range.Start, $blankPara.Range.Start)

$newSectionBody = @'
<w:p><w:r><w:t>{{#if HasFailedExpressions}}</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Failed Expressions</w:t></w:r></w:p><w:p><w:r><w:t>The following expressions could not be evaluated:</w:t></w:r></w:p><w:p/><w:tbl><w:tblPr><w:tblBorders><w:top w:val="single" w:sz="4"/><w:bottom w:val="single" w:sz="4"/><w:left w:val="single" w:sz="4"/><w:right w:val="single" w:sz="4"/><w:insideH w:val="single" w:sz="4"/><w:insideV w:val="single" w:sz="4"/></w:tblBorders></w:tblPr><w:tr><w:tc><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Expression</w:t></w:r></w:p></w:tc><w:tc><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Error</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:r><w:t>{{#foreach FailedExpressions}}</w:t></w:r></w:p></w:tc><w:tc><w:p><w:r><w:t/></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:r><w:t>{{VariableName}}</w:t></w:r></w:p></w:tc><w:tc><w:p><w:r><w:t>{{Message}}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:p><w:r><w:t>{{/foreach}}</w:t></w:r></w:p></w:tc><w:tc><w:p><w:r><w:t/></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:p><w:r><w:t>{{/if}}</w:t></w:r></w:p>
'@

$packageXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newSectionBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$insertBefore.InsertXML($packageXml)
